$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text changed from "Ready for handoff" to "In Translation" everywhere
# it appears (Overview summary columns E & F, and the per-language "Status"
# column C on the zh-cn / de-de report sheets).
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# Column widths shrink to re-fit the (shorter) new status text.
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()
$wsZhCn.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(3).AutoFit()
